# Generate Report for Handback
# Adds a new handback entry (b16978cc-c86d-4ccc-8566-95c3e8ba65e4) as row 4
# on the Overview, zh-cn and de-de sheets, and resizes their tables.

$wb = $excel.ActiveWorkbook

$guid = "b16978cc-c86d-4ccc-8566-95c3e8ba65e4"
$mdName = "$guid.md"
$mdDisplay = "e2e\$guid.md"
$statusInSync = "Handed back: in sync with en-US"

$hyperlinkColor = 15570276  # OLE/BGR form of RGB(100,149,237) == #FF6495ED

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $mdDisplay
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = "2016-10-21 04:40:12"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ceb469cef77ccb6416b2aafb9a125826abb923c8/e2e/$mdName", "", "", $mdDisplay)
$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = $hyperlinkColor

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> row 4
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf = "$guid.a5b39d683d79a4db47f9296e81376542f83da4b5.zh-cn.xlf"

$wsZhCn.Range("A4").Value = $mdName
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $statusInSync
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = "2016-10-21 04:40:00"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I4").Value = $mdName
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = "2016-10-21 04:40:54"
$wsZhCn.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ceb469cef77ccb6416b2aafb9a125826abb923c8/e2e/$mdName", "", "", $mdName)
$wsZhCn.Range("A4").Font.Underline = $true
$wsZhCn.Range("A4").Font.Color = $hyperlinkColor

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c3d4e1a9276ef64be3a7410c3b4fcb2ef2cf0dd7/e2e/$mdName", "", "", $mdName)
$wsZhCn.Range("I4").Font.Underline = $true
$wsZhCn.Range("I4").Font.Color = $hyperlinkColor

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> row 4
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf = "$guid.a5b39d683d79a4db47f9296e81376542f83da4b5.de-de.xlf"

$wsDeDe.Range("A4").Value = $mdName
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $statusInSync
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = "2016-10-21 04:40:12"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I4").Value = $mdName
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = "2016-10-21 04:41:13"
$wsDeDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ceb469cef77ccb6416b2aafb9a125826abb923c8/e2e/$mdName", "", "", $mdName)
$wsDeDe.Range("A4").Font.Underline = $true
$wsDeDe.Range("A4").Font.Color = $hyperlinkColor

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9e8b935f19f3ae14a4c2a1b0d7e2c5f6a4b1d7e0/e2e/$mdName", "", "", $mdName)
$wsDeDe.Range("I4").Font.Underline = $true
$wsDeDe.Range("I4").Font.Color = $hyperlinkColor

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))
